$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: hours logged on 2023-02-09 (F5 date)
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 4

# Row 10: hours logged on 2023-02-14 (F10 date)
$ws.Range("G10").Value = 2.5
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 2

# Recalculate totals / dependent formulas
$excel.Calculate()

# Update the active selection to mirror where the author ended up clicking
$null = $ws.Range("I17").Select()
